$wb = $excel.ActiveWorkbook

# --- "Core" sheet: the per-question "code" column (A) used to hold a
# dedicated testchartcorecode<N> string. The fixture now reuses the
# question name (column B) as the code value instead. ---
$core = $wb.Worksheets.Item("Core")
$core.Range("A2").Value = $core.Range("B2").Value2
$core.Range("A3").Value = $core.Range("B3").Value2
$core.Range("A4").Value = $core.Range("B4").Value2
$core.Range("A5").Value = $core.Range("B5").Value2

# --- "Test Chart" sheet: the first question's code changes from the
# generic "testchartcode0" placeholder to the real "PatientChartingDate"
# code. ---
$chart = $wb.Worksheets.Item("Test Chart")
$chart.Range("A2").Value = "PatientChartingDate"
$chart.Range("A2").ShrinkToFit = $false
